$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the price of Mango (row 2, column D) from 299 to 150
$ws.Range("D2").Value = 150
